# Refresh cryptocurrency price/volume snapshot (GitHub Actions data pull).
# Source cells are plain text (inline strings) in the workbook, e.g. "64.881.08"
# or "  -2.26%  ", so force text format before writing to stop Excel's COM layer
# from re-parsing number-shaped strings (e.g. "567.60" -> 567.6, dropping the
# trailing zero) or autocorrecting them in any other way.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '64.881.08'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -2.26%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.157.39'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -7.68%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '567.60'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -2.90%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '170.88'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -4.96%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.610'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -1.74%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '3.157.04'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -7.58%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.125'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -5.46%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.58'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -5.52%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -4.66%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.703.38'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -7.72%  '
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +1.02%  '
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -6.78%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.838.72'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -2.26%  '
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -5.59%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.159.46'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -7.47%  '
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -3.14%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.86'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -7.04%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '356.41'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -3.44%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.26'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -4.45%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.01%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '68.86'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -6.06%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -6.63%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -7.46%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.67'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -1.56%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -1.95%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +0.11%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.13%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -4.27%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -7.51%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '21.99'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -6.05%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -5.57%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.21'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -6.05%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '154.16'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -5.99%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -3.91%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '26.11'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -5.50%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -2.46%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -2.89%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.669.42'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -2.07%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.18'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -5.51%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '6.02'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -4.66%  '
$ws.Range('B45').NumberFormat = '@'
$ws.Range('B45').Value = 'OKB'
$ws.Range('C45').NumberFormat = '@'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '39.37'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -1.42%  '
$ws.Range('B46').NumberFormat = '@'
$ws.Range('B46').Value = 'Hedera'
$ws.Range('C46').NumberFormat = '@'
$ws.Range('C46').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0659'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -4.18%  '
$ws.Range('B47').NumberFormat = '@'
$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').NumberFormat = '@'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '24.29'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -3.35%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '328.71'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -2.37%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0275'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -4.26%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -2.27%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.03%  '
